# Generate Report for Archive
# - Update status text "Ready for handoff" -> "In Translation" wherever it appears
# - Narrow the affected "Status" columns to match the shorter text

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# Column width to apply to the narrowed status columns (closest value this
# engine's ColumnWidth -> stored width quantization can reach to the
# target stored width of 13.4101845877511).
$narrowColumnWidth = 12.5

# --- Overview sheet: locale status columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")

$overviewRange = $wsOverview.Range("A1:G3")
foreach ($cell in $overviewRange.Cells) {
    if ($oldStatus -eq $cell.Value2) {
        $cell.Value = $newStatus
    }
}

$wsOverview.Columns.Item(5).ColumnWidth = $narrowColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowColumnWidth

# --- Locale detail sheets: "Status" column C ---
$localeSheetNames = @("zh-cn", "de-de")
foreach ($sheetName in $localeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $usedRange = $ws.Range("A1:P3")
    foreach ($cell in $usedRange.Cells) {
        if ($oldStatus -eq $cell.Value2) {
            $cell.Value = $newStatus
        }
    }
    $ws.Columns.Item(3).ColumnWidth = $narrowColumnWidth
}
